$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E608").Value = 10736547729
$ws.Range("F608").Value = 106653933643
$ws.Range("E609").Value = 10736804021
$ws.Range("F609").Value = 106653773500
$ws.Range("E610").Value = 10737123400
$ws.Range("F610").Value = 106653586309
$ws.Range("E611").Value = 10737390755
$ws.Range("F611").Value = 106653477432
$ws.Range("E612").Value = 10737701841
$ws.Range("F612").Value = 106653466308
$ws.Range("E613").Value = 10737967314
$ws.Range("F613").Value = 106653450259
$ws.Range("E614").Value = 10737036733
$ws.Range("F614").Value = 106653604844
$ws.Range("E615").Value = 10736979884
$ws.Range("F615").Value = 106653532127
$ws.Range("E616").Value = 10736395654
$ws.Range("F616").Value = 106653550629
$ws.Range("E617").Value = 10736089608
$ws.Range("F617").Value = 106653564502
$ws.Range("E618").Value = 10735950312
$ws.Range("F618").Value = 106653567934
$ws.Range("E619").Value = 10735679599
$ws.Range("F619").Value = 106653599103
$ws.Range("E620").Value = 10735571406
$ws.Range("F620").Value = 106653613333
$ws.Range("E621").Value = 10735382214
$ws.Range("F621").Value = 106653637485
$ws.Range("E622").Value = 10735200117
$ws.Range("F622").Value = 106653669136
$ws.Range("E623").Value = 10735859999
$ws.Range("F623").Value = 106653483333
$ws.Range("E624").Value = 10735936216
$ws.Range("F624").Value = 106653313266
$ws.Range("E625").Value = 10735924439
$ws.Range("F625").Value = 106653109593
$ws.Range("E626").Value = 10736308333
$ws.Range("F626").Value = 106653250000
$ws.Range("E627").Value = 10736198173
$ws.Range("F627").Value = 106657899109
$ws.Range("E628").Value = 10736191841
$ws.Range("F628").Value = 106657691051
$ws.Range("E629").Value = 10736147375
$ws.Range("F629").Value = 106657501917
$ws.Range("E630").Value = 10736143080
$ws.Range("F630").Value = 106657298731
$ws.Range("E631").Value = 10736147385
$ws.Range("F631").Value = 106657106741
$ws.Range("E632").Value = 10736171062
$ws.Range("F632").Value = 106656968129
$ws.Range("E633").Value = 10736160120
$ws.Range("F633").Value = 106656770990
$ws.Range("E634").Value = 10735892861
$ws.Range("F634").Value = 106657693403
$ws.Range("E635").Value = 10735646134
$ws.Range("F635").Value = 106657691604
$ws.Range("E636").Value = 10735428935
$ws.Range("F636").Value = 106657721250
$ws.Range("E637").Value = 10735425966
$ws.Range("F637").Value = 106657445006
$ws.Range("E638").Value = 10735379041
$ws.Range("F638").Value = 106657380339
$ws.Range("E639").Value = 10735425379
$ws.Range("F639").Value = 106657192396
$ws.Range("E640").Value = 10735441891
$ws.Range("F640").Value = 106657104415
$ws.Range("E641").Value = 10735442347
$ws.Range("F641").Value = 106657018311
$ws.Range("E642").Value = 10735420052
$ws.Range("F642").Value = 106656767670
$ws.Range("E643").Value = 10735411917
$ws.Range("F643").Value = 106656588103
$ws.Range("E644").Value = 10733340000
$ws.Range("F644").Value = 106656128333
$ws.Range("E645").Value = 10735108333
$ws.Range("F645").Value = 106656208334
$ws.Range("E646").Value = 10735108333
$ws.Range("F646").Value = 106656208334
$ws.Range("E647").Value = 10733062585
$ws.Range("F647").Value = 106656218246
$ws.Range("E648").Value = 10733075679
$ws.Range("F648").Value = 106656464724
$ws.Range("E649").Value = 10733086270
$ws.Range("F649").Value = 106656643530
$ws.Range("E650").Value = 10733100765
$ws.Range("F650").Value = 106656924248
$ws.Range("E651").Value = 10733116314
$ws.Range("F651").Value = 106656994887
$ws.Range("E652").Value = 10733100213
$ws.Range("F652").Value = 106657152686
$ws.Range("E653").Value = 10733108626
$ws.Range("F653").Value = 106657318456
$ws.Range("E654").Value = 10732938935
$ws.Range("F654").Value = 106656990879
$ws.Range("E655").Value = 10732820685
$ws.Range("F655").Value = 106656958753
$ws.Range("E656").Value = 10732611446
$ws.Range("F656").Value = 106656964438
$ws.Range("E657").Value = 10737903333
$ws.Range("F657").Value = 106671090000
$ws.Range("E658").Value = 10737863261
$ws.Range("F658").Value = 106671076288
$ws.Range("E659").Value = 10737977733
$ws.Range("F659").Value = 106671357678
$ws.Range("E660").Value = 10738008952
$ws.Range("F660").Value = 106671474803
$ws.Range("E661").Value = 10738136791
$ws.Range("F661").Value = 106671794054
$ws.Range("E662").Value = 10738272834
$ws.Range("F662").Value = 106672155612
$ws.Range("E663").Value = 10738424635
$ws.Range("F663").Value = 106672552374

$ws.Range("K610").Select()
